$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 6: new K-Fold result columns (2nd MLP cross-validation run) ---
# B6 used to hold the "(Still running)" placeholder string; it now holds the
# completed average-accuracy value for this run.
$ws.Range("B6").Value = 0.99971586139201996

# C6/E6/G6/I6: best-parameter dictionaries for the four additional folds,
# written left-to-right so the shared-string table picks them up in order.
$ws.Range("C6").Value = "{'clf__hidden_layer_sizes': (15, 15), 'hstgm__num_buckets': 4}"
$ws.Range("C6").NumberFormat = "0.0000000"
$ws.Range("C6").WrapText = $true

$ws.Range("D6").Value = 0.99965128443566198

$ws.Range("E6").Value = "{'clf__hidden_layer_sizes': (13, 13), 'hstgm__num_buckets': 15}"
$ws.Range("E6").NumberFormat = "0.0000000"
$ws.Range("E6").WrapText = $true

$ws.Range("F6").Value = 0.99971586139201996

$ws.Range("G6").Value = "{'clf__hidden_layer_sizes': (19, 19), 'hstgm__num_buckets': 6}"
$ws.Range("G6").NumberFormat = "0.0000000"
$ws.Range("G6").WrapText = $true

$ws.Range("H6").Value = 0.999638364373724

$ws.Range("I6").Value = "{'clf__hidden_layer_sizes': (20, 20), 'hstgm__num_buckets': 2}"
$ws.Range("I6").WrapText = $true

# J6 takes over the "(Still running)" note that used to live in B6.
$ws.Range("J6").Value = "(Still running)"

# K6 stays empty but picks up the row's wrap-text style.
$ws.Range("K6").WrapText = $true

# --- Page setup: scale the printout down and switch to landscape ---
$ws.PageSetup.LeftMargin = $excel.InchesToPoints(0.25)
$ws.PageSetup.RightMargin = $excel.InchesToPoints(0.25)
$ws.PageSetup.Zoom = 53
$ws.PageSetup.Orientation = 2

# --- Selection moves to K6 (last touched cell) ---
$ws.Range("K6").Select()
